$d = $word.ActiveDocument

function Replace-Text($find, $repl) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $repl, 2)
    if (-not $ok) {
        Write-Output "WARN: not found -> $find"
    }
}

# --- Paragraph 2 (header): "На квалификационную работу бакалавра Онюшева Артема Андреевича"
# Merge runs so the spellcheck markers around "Онюшева" disappear (text itself is unchanged).
Replace-Text "бакалавра Онюшева Артема Андреевича" "бакалавра Онюшева Артема Андреевича"

# --- Paragraph 3 (header): "На тему: «»" -> full title inserted between the guillemets
Replace-Text "На тему: «»" "На тему: «Применение методов машинного обучения (ML) для решения задач технического анализа при управлении активами на фондовом рынке»"

# --- Paragraph 4: "Квалификационная работа студента Онюшева А.А. посвящена."
# Merge runs so spellcheck markers around "Онюшева" disappear, and append the new description.
Replace-Text "студента Онюшева А.А. посвящена." "студента Онюшева А.А. посвящена разработке различных архитектур нейронных сетей на Python, используя модуль PyTorch, и проведении на них исследований эффективности применения методов машинного обучения для решения задач технического анализа при управлении активами на фондовом рынке."

# --- Paragraph 5: big "В работе рассмотрено. ... Результатом ... стало." rewrite
Replace-Text "В работе рассмотрено. " "В работе рассмотрены три различных архитектуры: MLP, CNN, Transformer. "
Replace-Text "Также реализована. " "Также реализован метод дообучения нейронных сетей. "
Replace-Text "Разработаны." "Приобретены навыки работы с модулем PyTorch."
Replace-Text "Освоена работа. Приобретены навыки работы с. " "Проведено множество исследований. "
Replace-Text "Результатом практической деятельности студента в рамках данной работы стало." "Результат практической деятельности студента в рамках данной работы демонстрирует высокий уровень владения инструментами модуля PyTorch, понимание нюансов различных архитектур нейронных сетей, а также глубокое понимание важности использования различных оценочных метрик в задачах такого типа."

# --- Paragraph 7: "К достоинствам работы следует отнести. Кроме того, нужно отметить."
Replace-Text "К достоинствам работы следует отнести. " "К достоинствам работы следует отнести высокий уровень оптимизации и удобства использования. "
Replace-Text "Кроме того, нужно отметить." "Кроме того, нужно отметить хорошую модульность кода. Практическая значимость работы заключается в подготовке необходимых исследований для дальнейших изучения применения методом машинного обучения для решения задач технического анализа при управлении активами на фондовом рынке."

# --- Paragraph 8: "К замечаниям ... подробного описания классов разработанных систем."
Replace-Text "подробного описания классов разработанных систем." "подробной документации написанных функций и классов."

# --- Last signature paragraph: merge "», " and "к.т.н" runs (drops the spellcheck markers around "к.т.н")
Replace-Text "Афанасьева», к.т.н" "Афанасьева», к.т.н"

Write-Output "All replacements attempted."
